$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 320, shifting existing rows 320:426 down to 321:427.
$ws.Rows.Item(320).Insert()

# Populate the newly inserted row 320 with the new record's data.
$ws.Range("A320").Value = 5
$ws.Range("B320").Value = "Macroferia Regional de Talca"
$ws.Range("C320").Value = "Maule"
$ws.Range("D320").Value = 44588
$ws.Range("E320").Value = 7
$ws.Range("F320").Value = "Fruta"
$ws.Range("G320").Value = 100102
$ws.Range("H320").Value = "Cítricos"
$ws.Range("I320").Value = 100102005
$ws.Range("J320").Value = "Naranja"
$ws.Range("K320").Value = "Valencia"
$ws.Range("L320").Value = "Primera"
$ws.Range("M320").Value = 300
$ws.Range("N320").Value = 9000
$ws.Range("O320").Value = 9000
$ws.Range("P320").Value = 9000
$ws.Range("Q320").Value = "$/bandeja 15 kilos granel"
$ws.Range("R320").Value = "Región de O'Higgins"
$ws.Range("S320").Value = 600
$ws.Range("T320").Value = 15
